$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in the title row
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Abril de 2020 a las 00:20"

# Update country rows: names (reordered country list) and statistics
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 240660
$ws.Cells.Item(4, 3).Value = 25657
$ws.Cells.Item(4, 4).Value = 10400
$ws.Cells.Item(4, 5).Value = 224449
$ws.Cells.Item(4, 6).Value = 5421
$ws.Cells.Item(4, 7).Value = 709
$ws.Cells.Item(4, 8).Value = 5811

$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Cells.Item(7, 2).Value = 84794
$ws.Cells.Item(7, 3).Value = 6813
$ws.Cells.Item(7, 4).Value = 22440
$ws.Cells.Item(7, 5).Value = 61247
$ws.Cells.Item(7, 6).Value = 3936
$ws.Cells.Item(7, 7).Value = 176
$ws.Cells.Item(7, 8).Value = 1107

$ws.Cells.Item(16, 1).Value = "Canada"
$ws.Cells.Item(16, 2).Value = 11283
$ws.Cells.Item(16, 3).Value = 1552
$ws.Cells.Item(16, 4).Value = 1979
$ws.Cells.Item(16, 5).Value = 9131
$ws.Cells.Item(16, 6).Value = 120
$ws.Cells.Item(16, 7).Value = 59
$ws.Cells.Item(16, 8).Value = 173

$ws.Cells.Item(23, 1).Value = "Australia"
$ws.Cells.Item(23, 2).Value = 5279
$ws.Cells.Item(23, 3).Value = 231
$ws.Cells.Item(23, 4).Value = 585
$ws.Cells.Item(23, 5).Value = 4666
$ws.Cells.Item(23, 6).Value = 50
$ws.Cells.Item(23, 7).Value = 5
$ws.Cells.Item(23, 8).Value = 28

$ws.Cells.Item(35, 1).Value = "India"
$ws.Cells.Item(35, 2).Value = 2543
$ws.Cells.Item(35, 3).Value = 545
$ws.Cells.Item(35, 4).Value = 191
$ws.Cells.Item(35, 5).Value = 2280
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 14
$ws.Cells.Item(35, 8).Value = 72

$ws.Cells.Item(134, 1).Value = "Jamaica"
$ws.Cells.Item(134, 2).Value = 47
$ws.Cells.Item(134, 3).Value = 3
$ws.Cells.Item(134, 4).Value = 2
$ws.Cells.Item(134, 5).Value = 42
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 3

$ws.Cells.Item(135, 1).Value = "Guatemala"
$ws.Cells.Item(135, 2).Value = 47
$ws.Cells.Item(135, 3).Value = 8
$ws.Cells.Item(135, 4).Value = 12
$ws.Cells.Item(135, 5).Value = 34
$ws.Cells.Item(135, 6).Value = 1
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 1

$ws.Cells.Item(136, 1).Value = "Barbados"
$ws.Cells.Item(136, 2).Value = 46
$ws.Cells.Item(136, 3).Value = 1
$ws.Cells.Item(136, 4).Value = 0
$ws.Cells.Item(136, 5).Value = 46
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 0

$ws.Cells.Item(137, 1).Value = "Uganda"
$ws.Cells.Item(137, 2).Value = 45
$ws.Cells.Item(137, 3).Value = 1
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 5).Value = 45
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0

$ws.Cells.Item(167, 1).Value = "Benin"
$ws.Cells.Item(167, 2).Value = 13
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 1
$ws.Cells.Item(167, 5).Value = 12
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

$ws.Cells.Item(168, 1).Value = "Santa Lucia"
$ws.Cells.Item(168, 2).Value = 13
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 1
$ws.Cells.Item(168, 5).Value = 12
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

$ws.Cells.Item(172, 1).Value = "Granada"
$ws.Cells.Item(172, 2).Value = 10
$ws.Cells.Item(172, 3).Value = 1
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 10
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(174, 1).Value = "Mozambique"
$ws.Cells.Item(174, 2).Value = 10
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 10
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

$ws.Cells.Item(175, 1).Value = "Surinam"
$ws.Cells.Item(175, 2).Value = 10
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 10
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

$ws.Cells.Item(176, 1).Value = "Laos"
$ws.Cells.Item(176, 2).Value = 10
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 10
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(181, 2).Value = 9
$ws.Cells.Item(181, 3).Value = 1
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 9
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

$ws.Cells.Item(182, 1).Value = "Zimbabue"
$ws.Cells.Item(182, 2).Value = 9
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 8
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 1

$ws.Cells.Item(183, 1).Value = "Montserrat"
$ws.Cells.Item(183, 2).Value = 9
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 7
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 2

$ws.Cells.Item(187, 1).Value = "Fiyi"
$ws.Cells.Item(187, 2).Value = 7
$ws.Cells.Item(187, 3).Value = 2
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 7
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

$ws.Cells.Item(188, 1).Value = "Santa Sede"
$ws.Cells.Item(188, 2).Value = 7
$ws.Cells.Item(188, 3).Value = 1
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

$ws.Cells.Item(190, 1).Value = "Nepal"
$ws.Cells.Item(190, 2).Value = 6
$ws.Cells.Item(190, 3).Value = 1
$ws.Cells.Item(190, 4).Value = 1
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

$ws.Cells.Item(191, 1).Value = "Cabo Verde"
$ws.Cells.Item(191, 2).Value = 6
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 1

$ws.Cells.Item(195, 1).Value = "Somalia"
$ws.Cells.Item(195, 2).Value = 5
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 1
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

$ws.Cells.Item(196, 1).Value = "Nicaragua"
$ws.Cells.Item(196, 2).Value = 5
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 4
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 1

$ws.Cells.Item(200, 1).Value = "Malaui"
$ws.Cells.Item(200, 2).Value = 3
$ws.Cells.Item(200, 3).Value = 3
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 3
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0

$ws.Cells.Item(201, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(201, 2).Value = 3
$ws.Cells.Item(201, 3).Value = 0
$ws.Cells.Item(201, 4).Value = 0
$ws.Cells.Item(201, 5).Value = 3
$ws.Cells.Item(201, 6).Value = 0
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 0

$ws.Cells.Item(204, 1).Value = "Anguila"
$ws.Cells.Item(204, 2).Value = 3
$ws.Cells.Item(204, 3).Value = 1
$ws.Cells.Item(204, 4).Value = 0
$ws.Cells.Item(204, 5).Value = 3
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 0

$ws.Cells.Item(205, 1).Value = "Burundi"
$ws.Cells.Item(205, 2).Value = 3
$ws.Cells.Item(205, 3).Value = 1
$ws.Cells.Item(205, 4).Value = 0
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

